$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 452
$ws.Cells.Item(80, 9).Value = 364.8
$ws.Cells.Item(80, 10).Value = 888
$ws.Cells.Item(80, 11).Value = 1094.4
$ws.Cells.Item(80, 12).Value = 2664
$ws.Cells.Item(80, 13).Value = -96.40000000000009
$ws.Cells.Item(80, 14).Value = -4660

$ws.Cells.Item(83, 8).Value = 452
$ws.Cells.Item(83, 9).Value = 364.8
$ws.Cells.Item(83, 10).Value = 888
$ws.Cells.Item(83, 11).Value = 3283.2
$ws.Cells.Item(83, 12).Value = 7992
$ws.Cells.Item(83, 13).Value = 1708.8
$ws.Cells.Item(83, 14).Value = -17976

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 338.16666
$ws.Cells.Item(2, 9).Value = 294.14285
$ws.Cells.Item(2, 11).Value = 294.14285
$ws.Cells.Item(2, 13).Value = -181.14285

$ws.Cells.Item(61, 8).Value = 1459
$ws.Cells.Item(61, 9).Value = 873.75
$ws.Cells.Item(61, 11).Value = 873.75
$ws.Cells.Item(61, 13).Value = -661.75

$ws.Cells.Item(76, 8).Value = 28548.25
$ws.Cells.Item(76, 10).Value = 28548.25
$ws.Cells.Item(76, 12).Value = 28548.25
$ws.Cells.Item(76, 14).Value = -29224.25

$ws.Cells.Item(79, 8).Value = 28548.25
$ws.Cells.Item(79, 10).Value = 28548.25
$ws.Cells.Item(79, 12).Value = 28548.25
$ws.Cells.Item(79, 14).Value = -30888.25

$ws.Cells.Item(116, 8).Value = 338.16666
$ws.Cells.Item(116, 9).Value = 294.14285
$ws.Cells.Item(116, 11).Value = 294.14285
$ws.Cells.Item(116, 13).Value = 1999.85715

$ws.Cells.Item(136, 8).Value = 1459
$ws.Cells.Item(136, 9).Value = 873.75
$ws.Cells.Item(136, 11).Value = 2621.25
$ws.Cells.Item(136, 13).Value = -71.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 338.16666
$ws.Cells.Item(3, 9).Value = 294.14285
$ws.Cells.Item(3, 11).Value = 294.14285
$ws.Cells.Item(3, 13).Value = -180.14285

$ws.Cells.Item(7, 8).Value = 156.33333
$ws.Cells.Item(7, 9).Value = 203
$ws.Cells.Item(7, 10).Value = 133
$ws.Cells.Item(7, 11).Value = 203
$ws.Cells.Item(7, 12).Value = 133
$ws.Cells.Item(7, 13).Value = -90
$ws.Cells.Item(7, 14).Value = -359

$ws.Cells.Item(107, 8).Value = 45142.445
$ws.Cells.Item(107, 9).Value = 67132.336
$ws.Cells.Item(107, 11).Value = 67132.336
$ws.Cells.Item(107, 13).Value = -65212.336

$ws.Cells.Item(134, 8).Value = 904
$ws.Cells.Item(134, 9).Value = 904
$ws.Cells.Item(134, 11).Value = 2712
$ws.Cells.Item(134, 13).Value = -177

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 459
$ws.Cells.Item(16, 9).Value = 365.66666
$ws.Cells.Item(16, 10).Value = 599
$ws.Cells.Item(16, 11).Value = 365.66666
$ws.Cells.Item(16, 12).Value = 599
$ws.Cells.Item(16, 13).Value = -78.66665999999998
$ws.Cells.Item(16, 14).Value = -1173

$ws.Cells.Item(58, 8).Value = 2011.5
$ws.Cells.Item(58, 9).Value = 2011.5
$ws.Cells.Item(58, 11).Value = 2011.5
$ws.Cells.Item(58, 13).Value = -1808.5

$ws.Cells.Item(70, 8).Value = 50000
$ws.Cells.Item(70, 9).Value = 50000
$ws.Cells.Item(70, 11).Value = 50000
$ws.Cells.Item(70, 13).Value = -49685

$ws.Cells.Item(73, 8).Value = 50000
$ws.Cells.Item(73, 9).Value = 50000
$ws.Cells.Item(73, 11).Value = 50000
$ws.Cells.Item(73, 13).Value = -48908

$ws.Cells.Item(74, 8).Value = 82656
$ws.Cells.Item(74, 10).Value = 82656
$ws.Cells.Item(74, 12).Value = 82656
$ws.Cells.Item(74, 14).Value = -84404

$ws.Cells.Item(77, 8).Value = 82656
$ws.Cells.Item(77, 10).Value = 82656
$ws.Cells.Item(77, 12).Value = 247968
$ws.Cells.Item(77, 14).Value = -256704

$ws.Cells.Item(88, 8).Value = 44468.4
$ws.Cells.Item(88, 10).Value = 44468.4
$ws.Cells.Item(88, 12).Value = 44468.4
$ws.Cells.Item(88, 14).Value = -45280.4

$ws.Cells.Item(91, 8).Value = 44468.4
$ws.Cells.Item(91, 10).Value = 44468.4
$ws.Cells.Item(91, 12).Value = 44468.4
$ws.Cells.Item(91, 14).Value = -47276.4

$ws.Cells.Item(113, 8).Value = 459
$ws.Cells.Item(113, 9).Value = 365.66666
$ws.Cells.Item(113, 10).Value = 599
$ws.Cells.Item(113, 11).Value = 365.66666
$ws.Cells.Item(113, 12).Value = 599
$ws.Cells.Item(113, 13).Value = 1804.33334
$ws.Cells.Item(113, 14).Value = -4939

$ws.Cells.Item(132, 8).Value = 1179.8
$ws.Cells.Item(132, 9).Value = 974.75
$ws.Cells.Item(132, 11).Value = 2924.25
$ws.Cells.Item(132, 13).Value = -394.25

$ws.Cells.Item(136, 8).Value = 2011.5
$ws.Cells.Item(136, 9).Value = 2011.5
$ws.Cells.Item(136, 11).Value = 6034.5
$ws.Cells.Item(136, 13).Value = -3484.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).ClearContents()
$ws.Cells.Item(69, 14).ClearContents()

$ws.Cells.Item(72, 8).Value = 0
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 13).ClearContents()
$ws.Cells.Item(72, 14).ClearContents()

$ws.Cells.Item(121, 8).Value = 497
$ws.Cells.Item(121, 9).Value = 208.8
$ws.Cells.Item(121, 10).Value = 977.3333
$ws.Cells.Item(121, 11).Value = 626.4000000000001
$ws.Cells.Item(121, 12).Value = 2931.9999
$ws.Cells.Item(121, 13).Value = 683.5999999999999
$ws.Cells.Item(121, 14).Value = -5551.9999

$ws.Cells.Item(131, 8).Value = 1799.5385
$ws.Cells.Item(131, 9).Value = 812.375
$ws.Cells.Item(131, 11).Value = 2437.125
$ws.Cells.Item(131, 13).Value = 2602.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(98, 8).Value = 8503.200000000001
$ws.Cells.Item(98, 10).Value = 8503.200000000001
$ws.Cells.Item(98, 12).Value = 8503.200000000001
$ws.Cells.Item(98, 14).Value = -14493.2

$ws.Cells.Item(113, 8).Value = 922.25
$ws.Cells.Item(113, 9).Value = 922.25
$ws.Cells.Item(113, 11).Value = 922.25
$ws.Cells.Item(113, 13).Value = 1247.75

$ws.Cells.Item(132, 8).Value = 1279.4
$ws.Cells.Item(132, 9).Value = 1279.4
$ws.Cells.Item(132, 11).Value = 3838.2
$ws.Cells.Item(132, 13).Value = -1308.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 206.57143
$ws.Cells.Item(22, 9).Value = 199.33333
$ws.Cells.Item(22, 10).Value = 250
$ws.Cells.Item(22, 11).Value = 199.33333
$ws.Cells.Item(22, 12).Value = 250
$ws.Cells.Item(22, 13).Value = 95.66667000000001
$ws.Cells.Item(22, 14).Value = -840

$ws.Cells.Item(27, 8).Value = 206.57143
$ws.Cells.Item(27, 9).Value = 199.33333
$ws.Cells.Item(27, 10).Value = 250
$ws.Cells.Item(27, 11).Value = 199.33333
$ws.Cells.Item(27, 12).Value = 250
$ws.Cells.Item(27, 13).Value = -92.33332999999999
$ws.Cells.Item(27, 14).Value = -464

$ws.Cells.Item(55, 8).Value = 1074.2273
$ws.Cells.Item(55, 9).Value = 842.53845
$ws.Cells.Item(55, 11).Value = 842.53845
$ws.Cells.Item(55, 13).Value = -669.53845

$ws.Cells.Item(61, 8).Value = 1900
$ws.Cells.Item(61, 9).Value = 1900
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 1900
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -1698
$ws.Cells.Item(61, 14).ClearContents()

$ws.Cells.Item(113, 8).Value = 1900
$ws.Cells.Item(113, 9).Value = 1900
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 1900
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = 270
$ws.Cells.Item(113, 14).ClearContents()

$ws.Cells.Item(134, 8).Value = 25000
$ws.Cells.Item(134, 10).Value = 25000
$ws.Cells.Item(134, 12).Value = 25000
$ws.Cells.Item(134, 14).Value = -35140

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64, 8).Value = 10526
$ws.Cells.Item(64, 10).Value = 10526
$ws.Cells.Item(64, 12).Value = 10526
$ws.Cells.Item(64, 14).Value = -11022

$ws.Cells.Item(67, 8).Value = 10526
$ws.Cells.Item(67, 10).Value = 10526
$ws.Cells.Item(67, 12).Value = 10526
$ws.Cells.Item(67, 14).Value = -12242

$ws.Cells.Item(76, 8).Value = 0
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 14).ClearContents()

$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 14).ClearContents()

$ws.Cells.Item(107, 8).Value = 590.2
$ws.Cells.Item(107, 9).Value = 425
$ws.Cells.Item(107, 11).Value = 1275
$ws.Cells.Item(107, 13).Value = 645

$ws.Cells.Item(113, 8).Value = 504.7143
$ws.Cells.Item(113, 9).Value = 290.5
$ws.Cells.Item(113, 11).Value = 871.5
$ws.Cells.Item(113, 13).Value = 1298.5
